$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - becomes the former row 4 match (EGYPT - PREMIER LEAGUE: ZED vs Al Ahly), new Id APE9ifU7
# (column B / Date is left untouched since it does not change)
$ws.Cells.Item(3, 1).Value = "APE9ifU7"
$row3 = @(
  "15:00","EGYPT - PREMIER LEAGUE","ZED","Al Ahly",
  5.5,3.8,1.57,5.3,2.18,2.15,1.06,7.5,1.28,3.35,1.85,1.9,1.39,2.75,1.87,1.83,
  14.5,32,17,100,55,55,7.5,7.4,16.5,80,600,6.6,7.2,8,11.25,13,27,6.9,30,35,
  200,200,450,2.75,7.8,75,3.4,7.7,18,25,60,250,51,51
)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 3).Value = $row3[$i]
}

# Row 4 - becomes a new match (ENGLAND - CHAMPIONSHIP: West Brom vs Burnley), new Id Emn9XsgJ
# (column B / Date is left untouched since it does not change)
$ws.Cells.Item(4, 1).Value = "Emn9XsgJ"
$row4 = @(
  "17:00","ENGLAND - CHAMPIONSHIP","West Brom","Burnley",
  2.38,3.2,3.1,3.2,2,3.75,1.08,8,1.44,2.75,2.38,1.57,1.5,2.5,2,1.73,
  6.5,11,10,23,21,34,7.5,6,17,51,401,8,15,12,34,29,41,4.33,15,26,
  51,81,201,2.5,8.5,67,5,19,29,67,101,251,126,151
)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 3).Value = $row4[$i]
}

# Row 5 - a few odds values updated
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("AN5").Value = 3.5
$ws.Range("AO5").Value = 7.5

# Row 6 - a few odds values updated
$ws.Range("G6").Value = 7.5
$ws.Range("J6").Value = 6.5
$ws.Range("L6").Value = 1.8
$ws.Range("N6").Value = 10
$ws.Range("Q6").Value = 1.62
$ws.Range("R6").Value = 2.25
$ws.Range("U6").Value = 1.91
$ws.Range("V6").Value = 1.8
$ws.Range("Y6").Value = 23
$ws.Range("AB6").Value = 51
$ws.Range("AH6").Value = 8
$ws.Range("AI6").Value = 7
$ws.Range("AK6").Value = 9
$ws.Range("AM6").Value = 26
$ws.Range("AP6").Value = 41
$ws.Range("AR6").Value = 151
$ws.Range("AU6").Value = 9
$ws.Range("AW6").Value = 3.5
